# Locate the "Summary" body paragraph (the one that currently starts with
# "This is a brief explanation...") and the empty paragraph directly after
# it, then replace that whole span with the rewritten Summary paragraph,
# a blank spacer paragraph, a new "Script" Heading2 paragraph, and the new
# Script body paragraph (which keeps the _GoBack bookmark at its end).
$d = $word.ActiveDocument

$summaryIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "This is a brief explanation*") {
        $summaryIndex = $i
        break
    }
}

$startPara = $d.Paragraphs($summaryIndex)
$endPara = $d.Paragraphs($summaryIndex + 1)
$target = $d.Range($startPara.Range.Start, $endPara.Range.End)

$newBodyXml = @'
<w:p><w:r><w:t xml:space="preserve">This is a brief explanation about how the Key Pressing Script functions. The script is used to have music notes interact with the keynotes. This is using a Raycast to check the distance between the </w:t></w:r><w:r><w:t xml:space="preserve">piano keys </w:t></w:r><w:r><w:t xml:space="preserve">and </w:t></w:r><w:r><w:t>the music note. The distance for the perfect hit has been set to 1 unit from the piano. Anything between 1 to 0.75 units would give the player points. Anything between 1 and 3 units would be considered good but not perfect, which will result in points still being given.</w:t></w:r><w:r><w:t xml:space="preserve"> Anything more than 3 units is terrible, nothing will be given.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Script</w:t></w:r></w:p><w:p><w:r><w:t>KeyPressing Script – The Raycast starts from the center of the piano key itself</w:t></w:r><w:r><w:t xml:space="preserve"> and shoots up in the dir</w:t></w:r><w:r><w:t>ection of on-coming music notes.</w:t></w:r><w:r><w:t xml:space="preserve"> The script checks </w:t></w:r><w:r><w:t xml:space="preserve">to see if the user hits the </w:t></w:r><w:r><w:t>music notes either too early, on spot or too late.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$pkgHeader = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pkg = $pkgHeader + $newBodyXml + $pkgFooter

[void]$target.InsertXML($pkg)

Write-Output "done"
